# daily auto push: 2026-01-19 22:35 UTC
# A new observation for 2026/01/20 (weekday 火) at time-slot 6 needs to be
# inserted into the log. It belongs right after the existing 2026/01/20
# row (row 674), so every row from the old 674 onward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 674..715 down to 675..716 (this also grows the used range /
# <dimension> from D715 to D716 automatically).
$ws.Rows(674).Insert()

# Fill in the newly opened row 674 with the inserted record.
# The date is entered with a leading apostrophe so Excel stores it as
# literal text ("2026/01/20") instead of auto-converting it to a date
# serial number, matching how the rest of column A is stored; resetting
# the style back to Normal afterwards drops the quote-prefix formatting
# that the apostrophe would otherwise leave behind.
$ws.Range("A674").Value = "'2026/01/20"
$ws.Range("A674").Style = "Normal"
$ws.Range("B674").Value = "火"
$ws.Range("C674").Value = 6
$ws.Range("D674").Value = 201
